$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$edits = @(
    @{Ref="D2"; Val="51.617.18"},
    @{Ref="E2"; Val="  +1.11%  "},
    @{Ref="D3"; Val="2.988.63"},
    @{Ref="E3"; Val="  +2.65%  "},
    @{Ref="E4"; Val="  +0.03%  "},
    @{Ref="D5"; Val="382.79"},
    @{Ref="E5"; Val="  +3.04%  "},
    @{Ref="D6"; Val="104.57"},
    @{Ref="E6"; Val="  +2.61%  "},
    @{Ref="D7"; Val="0.547"},
    @{Ref="E7"; Val="  +1.22%  "},
    @{Ref="E8"; Val="  -0.01%  "},
    @{Ref="D9"; Val="0.597"},
    @{Ref="E9"; Val="  +2.24%  "},
    @{Ref="D10"; Val="37.44"},
    @{Ref="E10"; Val="  +1.55%  "},
    @{Ref="E11"; Val="  +0.11%  "},
    @{Ref="D12"; Val="0.0847"},
    @{Ref="E12"; Val="  +1.62%  "},
    @{Ref="D13"; Val="3.453.09"},
    @{Ref="E13"; Val="  +2.45%  "},
    @{Ref="D14"; Val="18.46"},
    @{Ref="E14"; Val="  +0.94%  "},
    @{Ref="D15"; Val="7.55"},
    @{Ref="E15"; Val="  +2.55%  "},
    @{Ref="D16"; Val="2.986.32"},
    @{Ref="E16"; Val="  +2.49%  "},
    @{Ref="D17"; Val="0.980"},
    @{Ref="E17"; Val="  +6.39%  "},
    @{Ref="D18"; Val="51.539.69"},
    @{Ref="E18"; Val="  +1.09%  "},
    @{Ref="D19"; Val="3.32"},
    @{Ref="E19"; Val="  +3.13%  "},
    @{Ref="D20"; Val="7.47"},
    @{Ref="E20"; Val="  +4.10%  "},
    @{Ref="D21"; Val="12.98"},
    @{Ref="E21"; Val="  +0.75%  "},
    @{Ref="D22"; Val="0.0₃0967"},
    @{Ref="E22"; Val="  +2.81%  "},
    @{Ref="D23"; Val="69.02"},
    @{Ref="E23"; Val="  +1.53%  "},
    @{Ref="D24"; Val="263.39"},
    @{Ref="E24"; Val="  +1.92%  "},
    @{Ref="D25"; Val="2.93"},
    @{Ref="E25"; Val="  +9.60%  "},
    @{Ref="D26"; Val="8.32"},
    @{Ref="E26"; Val="  +18.03%  "},
    @{Ref="D27"; Val="7.76"},
    @{Ref="E27"; Val="  +24.90%  "},
    @{Ref="D28"; Val="0.117"},
    @{Ref="E28"; Val="  +15.50%  "},
    @{Ref="D29"; Val="0.171"},
    @{Ref="E29"; Val="  +2.52%  "},
    @{Ref="D30"; Val="26.06"},
    @{Ref="E30"; Val="  +1.98%  "},
    @{Ref="E31"; Val="  +0.04%  "},
    @{Ref="D32"; Val="9.92"},
    @{Ref="E32"; Val="  +0.74%  "},
    @{Ref="D33"; Val="34.99"},
    @{Ref="E33"; Val="  +2.65%  "},
    @{Ref="E34"; Val="  -0.67%  "},
    @{Ref="E35"; Val="  -1.91%  "},
    @{Ref="D36"; Val="0.0453"},
    @{Ref="E36"; Val="  +7.86%  "},
    @{Ref="E37"; Val="  -0.12%  "},
    @{Ref="E38"; Val="  +1.81%  "},
    @{Ref="D39"; Val="17.14"},
    @{Ref="E39"; Val="  +0.85%  "},
    @{Ref="D40"; Val="2.59"},
    @{Ref="E40"; Val="  +0.60%  "},
    @{Ref="D41"; Val="1.85"},
    @{Ref="E41"; Val="  +0.83%  "},
    @{Ref="E42"; Val="  +3.41%  "},
    @{Ref="D43"; Val="122.30"},
    @{Ref="E43"; Val="  +2.56%  "},
    @{Ref="D44"; Val="21.87"},
    @{Ref="E44"; Val="  +0.24%  "},
    @{Ref="D45"; Val="0.282"},
    @{Ref="E45"; Val="  +20.10%  "},
    @{Ref="E46"; Val="  -2.49%  "},
    @{Ref="E47"; Val="  +2.89%  "},
    @{Ref="D48"; Val="3.30"},
    @{Ref="E48"; Val="  +5.20%  "},
    @{Ref="D49"; Val="2.036.00"},
    @{Ref="E49"; Val="  +0.84%  "},
    @{Ref="D50"; Val="0.0333"},
    @{Ref="E50"; Val="  +8.11%  "},
    @{Ref="D51"; Val="58.38"},
    @{Ref="E51"; Val="  +3.56%  "}
)

foreach ($edit in $edits) {
    $cell = $ws.Range($edit.Ref)
    $cell.NumberFormat = "@"
    $cell.Value = $edit.Val
    $cell.ClearFormats()
}
